$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.427.71"
$ws.Range("E2").Value = "  +2.48%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.385.53"
$ws.Range("E3").Value = "  +1.90%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.46"
$ws.Range("E5").Value = "  +1.59%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "179.14"
$ws.Range("E6").Value = "  +2.40%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.08%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.596"
$ws.Range("E8").Value = "  +1.42%  "

# Row 9
$ws.Range("E9").Value = "  +5.03%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.587"
$ws.Range("E10").Value = "  +2.07%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "48.37"
$ws.Range("E11").Value = "  +6.63%  "

# Row 12
$ws.Range("E12").Value = "  +3.23%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "701.36"
$ws.Range("E13").Value = "  +6.82%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.945.51"
$ws.Range("E14").Value = "  +2.10%  "

# Row 15
$ws.Range("E15").Value = "  +1.78%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "69.459.89"
$ws.Range("E16").Value = "  +2.67%  "

# Row 17
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.120"
$ws.Range("E17").Value = "  +1.67%  "

# Row 18
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.381.40"
$ws.Range("E18").Value = "  +1.62%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.64"
$ws.Range("E19").Value = "  +1.70%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.32"
$ws.Range("E20").Value = "  +3.54%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.903"
$ws.Range("E21").Value = "  +2.03%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.53"
$ws.Range("E22").Value = "  +3.28%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.16"
$ws.Range("E23").Value = "  +1.35%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "101.47"
$ws.Range("E24").Value = "  +3.26%  "

# Row 25
$ws.Range("E25").Value = "  +2.56%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.72"
$ws.Range("E26").Value = "  +2.24%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.62"
$ws.Range("E27").Value = "  +4.11%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "33.44"
$ws.Range("E28").Value = "  +0.43%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.64"
$ws.Range("E29").Value = "  +2.82%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.05"
$ws.Range("E30").Value = "  -2.19%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.16"
$ws.Range("E31").Value = "  +2.11%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "553.65"
$ws.Range("E32").Value = "  -2.45%  "

# Row 33
$ws.Range("E33").Value = "  +1.89%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "58.36"
$ws.Range("E34").Value = "  +3.72%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.728.30"
$ws.Range("E35").Value = "  +1.78%  "

# Row 36
$ws.Range("E36").Value = "  +5.83%  "

# Row 37
$ws.Range("E37").Value = "  -0.04%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.145"
$ws.Range("E38").Value = "  +11.72%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "34.94"
$ws.Range("E39").Value = "  +1.37%  "

# Row 40
$ws.Range("E40").Value = "  +3.59%  "

# Row 41
$ws.Range("E41").Value = "  +1.47%  "

# Row 42
$ws.Range("E42").Value = "  +3.75%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.339"
$ws.Range("E43").Value = "  +2.08%  "

# Row 44
$ws.Range("E44").Value = "  +3.76%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.24"
$ws.Range("E45").Value = "  -3.48%  "

# Row 46
$ws.Range("E46").Value = "  +3.01%  "

# Row 47
$ws.Range("E47").Value = "  +1.70%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.00"
$ws.Range("E48").Value = "  -0.07%  "

# Row 49
$ws.Range("E49").Value = "  -1.27%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.31"
$ws.Range("E50").Value = "  +3.40%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.63"
$ws.Range("E51").Value = "  -2.39%  "
